$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 421
$ws.Range("J33").Value = 533
$ws.Range("L33").Value = 533
$ws.Range("N33").Value = -991
$ws.Range("H98").Value = 2213.3096
$ws.Range("I98").Value = 684.3714
$ws.Range("K98").Value = 684.3714
$ws.Range("M98").Value = 813.6286
$ws.Range("H111").Value = 1471.6666
$ws.Range("I111").Value = 1321.5
$ws.Range("J111").Value = 1621.8334
$ws.Range("K111").Value = 3964.5
$ws.Range("L111").Value = 4865.5002
$ws.Range("M111").Value = -897.5
$ws.Range("N111").Value = -10999.5002
$ws.Range("H116").Value = 462266.53
$ws.Range("I116").Value = 1253094
$ws.Range("J116").Value = 10365.143
$ws.Range("K116").Value = 1253094
$ws.Range("L116").Value = 10365.143
$ws.Range("M116").Value = -1249652
$ws.Range("N116").Value = -17249.143
$ws.Range("H122").Value = 2213.3096
$ws.Range("I122").Value = 684.3714
$ws.Range("K122").Value = 2053.1142
$ws.Range("M122").Value = 396.8858
$ws.Range("H125").Value = 2282.5334
$ws.Range("I125").Value = 1012.6667
$ws.Range("J125").Value = 2600
$ws.Range("K125").Value = 9114.0003
$ws.Range("L125").Value = 23400
$ws.Range("M125").Value = -6654.0003
$ws.Range("N125").Value = -28320
$ws.Range("H137").Value = 4258.522
$ws.Range("I137").Value = 2137.75
$ws.Range("K137").Value = 6413.25
$ws.Range("M137").Value = -3863.25
$ws.Range("H141").Value = 5104.86
$ws.Range("I141").Value = 5155.0625
$ws.Range("J141").Value = 3900
$ws.Range("K141").Value = 15465.1875
$ws.Range("L141").Value = 11700
$ws.Range("M141").Value = -10285.1875
$ws.Range("N141").Value = -22060
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5335.547
$ws.Range("I32").Value = 4209.449
$ws.Range("J32").Value = 9014.134
$ws.Range("K32").Value = 4209.449
$ws.Range("L32").Value = 9014.134
$ws.Range("M32").Value = -3922.449
$ws.Range("N32").Value = -9588.134
$ws.Range("H61").Value = 826.3148
$ws.Range("I61").Value = 671.26666
$ws.Range("J61").Value = 1601.5555
$ws.Range("K61").Value = 671.26666
$ws.Range("L61").Value = 1601.5555
$ws.Range("M61").Value = -459.26666
$ws.Range("N61").Value = -2025.5555
$ws.Range("H122").Value = 2071.3044
$ws.Range("I122").Value = 1244.7693
$ws.Range("J122").Value = 3145.8
$ws.Range("K122").Value = 3734.3079
$ws.Range("L122").Value = 9437.400000000001
$ws.Range("M122").Value = -1284.3079
$ws.Range("N122").Value = -14337.4
$ws.Range("H132").Value = 2532.6
$ws.Range("I132").Value = 1327.6666
$ws.Range("K132").Value = 3982.9998
$ws.Range("M132").Value = -1452.9998
$ws.Range("H136").Value = 826.3148
$ws.Range("I136").Value = 671.26666
$ws.Range("J136").Value = 1601.5555
$ws.Range("K136").Value = 2013.79998
$ws.Range("L136").Value = 4804.666499999999
$ws.Range("M136").Value = 536.20002
$ws.Range("N136").Value = -9904.666499999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1832.579
$ws.Range("I99").Value = 1207.4615
$ws.Range("K99").Value = 1207.4615
$ws.Range("M99").Value = 290.5385000000001
$ws.Range("H134").Value = 1376.8948
$ws.Range("I134").Value = 989.8627300000001
$ws.Range("J134").Value = 4666.6665
$ws.Range("K134").Value = 2969.58819
$ws.Range("L134").Value = 13999.9995
$ws.Range("M134").Value = -434.5881900000004
$ws.Range("N134").Value = -19069.9995
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10002103
$ws.Range("I31").Value = 1160.0571
$ws.Range("K31").Value = 1160.0571
$ws.Range("M31").Value = -865.0571
$ws.Range("H34").Value = 10002103
$ws.Range("I34").Value = 1160.0571
$ws.Range("K34").Value = 1160.0571
$ws.Range("M34").Value = -958.0571
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H132").Value = 1176.2097
$ws.Range("I132").Value = 680.55554
$ws.Range("J132").Value = 2488.2354
$ws.Range("K132").Value = 2041.66662
$ws.Range("L132").Value = 7464.706200000001
$ws.Range("M132").Value = 488.33338
$ws.Range("N132").Value = -12524.7062
$ws.Range("H134").Value = 1660.875
$ws.Range("I134").Value = 659.03845
$ws.Range("J134").Value = 3521.4285
$ws.Range("K134").Value = 1977.11535
$ws.Range("L134").Value = 10564.2855
$ws.Range("M134").Value = 557.88465
$ws.Range("N134").Value = -15634.2855
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1441.3721
$ws.Range("I5").Value = 438.85
$ws.Range("J5").Value = 2313.1304
$ws.Range("K5").Value = 1316.55
$ws.Range("L5").Value = 6939.3912
$ws.Range("M5").Value = -1204.55
$ws.Range("N5").Value = -7163.3912
$ws.Range("H34").Value = 12344.889
$ws.Range("I34").Value = 36723.332
$ws.Range("J34").Value = 7469.2
$ws.Range("K34").Value = 110169.996
$ws.Range("L34").Value = 22407.6
$ws.Range("M34").Value = -110085.996
$ws.Range("N34").Value = -22575.6
$ws.Range("H56").Value = 6295.364
$ws.Range("I56").Value = 6295.364
$ws.Range("K56").Value = 6295.364
$ws.Range("M56").Value = -5765.364
$ws.Range("H113").Value = 603.6
$ws.Range("I113").Value = 562
$ws.Range("J113").Value = 666
$ws.Range("K113").Value = 1686
$ws.Range("L113").Value = 1998
$ws.Range("M113").Value = 484
$ws.Range("N113").Value = -6338
$ws.Range("H131").Value = 6494318
$ws.Range("J131").Value = 853.9859
$ws.Range("L131").Value = 2561.9577
$ws.Range("N131").Value = -12641.9577
$ws.Range("H132").Value = 2170.6562
$ws.Range("I132").Value = 634.1111
$ws.Range("J132").Value = 2771.913
$ws.Range("K132").Value = 5706.9999
$ws.Range("L132").Value = 24947.217
$ws.Range("M132").Value = -3176.9999
$ws.Range("N132").Value = -30007.217
$ws.Range("H134").Value = 4427.8965
$ws.Range("I134").Value = 3944.2144
$ws.Range("J134").Value = 4879.3335
$ws.Range("K134").Value = 11832.6432
$ws.Range("L134").Value = 14638.0005
$ws.Range("M134").Value = -6762.643199999999
$ws.Range("N134").Value = -24778.0005
$ws.Range("H135").Value = 1441.3721
$ws.Range("I135").Value = 438.85
$ws.Range("J135").Value = 2313.1304
$ws.Range("K135").Value = 3949.65
$ws.Range("L135").Value = 20818.1736
$ws.Range("M135").Value = -1414.65
$ws.Range("N135").Value = -25888.1736
$ws.Range("H140").Value = 2581.48
$ws.Range("I140").Value = 2761
$ws.Range("K140").Value = 8283
$ws.Range("M140").Value = -3103
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2194.19
$ws.Range("I126").Value = 2207.402
$ws.Range("J126").Value = 1767
$ws.Range("K126").Value = 6622.206
$ws.Range("L126").Value = 5301
$ws.Range("M126").Value = -4152.206
$ws.Range("N126").Value = -10241
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 39500
$ws.Range("I74").Value = 10000
$ws.Range("J74").Value = 46875
$ws.Range("K74").Value = 10000
$ws.Range("L74").Value = 46875
$ws.Range("M74").Value = -9002
$ws.Range("N74").Value = -48871
$ws.Range("H77").Value = 39500
$ws.Range("I77").Value = 10000
$ws.Range("J77").Value = 46875
$ws.Range("K77").Value = 30000
$ws.Range("L77").Value = 140625
$ws.Range("M77").Value = -25008
$ws.Range("N77").Value = -150609
$ws.Range("H132").Value = 11775.723
$ws.Range("I132").Value = 14354.368
$ws.Range("J132").Value = 8893.706
$ws.Range("K132").Value = 43063.104
$ws.Range("L132").Value = 26681.118
$ws.Range("M132").Value = -40533.104
$ws.Range("N132").Value = -31741.118
$ws.Range("H136").Value = 2828.8147
$ws.Range("I136").Value = 1465.6842
$ws.Range("J136").Value = 6066.25
$ws.Range("K136").Value = 4397.0526
$ws.Range("L136").Value = 18198.75
$ws.Range("M136").Value = -1847.0526
$ws.Range("N136").Value = -23298.75
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 496.53333
$ws.Range("I113").Value = 400
$ws.Range("J113").Value = 689.6
$ws.Range("K113").Value = 1200
$ws.Range("L113").Value = 2068.8
$ws.Range("M113").Value = 970
$ws.Range("N113").Value = -6408.8
$ws.Range("H132").Value = 5557278.5
$ws.Range("I132").Value = 1215.1163
$ws.Range("J132").Value = 19610850
$ws.Range("K132").Value = 3645.3489
$ws.Range("L132").Value = 58832550
$ws.Range("M132").Value = -1115.3489
$ws.Range("N132").Value = -58837610
$ws.Range("H136").Value = 2358.0667
$ws.Range("I136").Value = 512.5
$ws.Range("J136").Value = 8062.5454
$ws.Range("K136").Value = 1537.5
$ws.Range("L136").Value = 24187.6362
$ws.Range("M136").Value = 1012.5
$ws.Range("N136").Value = -29287.6362
